# Applies the "Modal eliminar, reubicacion de boton 'Nueva Tarea'" edit.
#
# Summary of the change (see commit message / xml diff):
#   - The last bullet ("Si agregas nuevas funcionalidades...") keeps its
#     <w:proofErr w:type="gramStart"/> ... <w:proofErr w:type="gramEnd"/>
#     pair tight around the run; the "_GoBack" bookmark that used to sit
#     at the end of that paragraph is relocated to the end of the new
#     "Actividad “favorito/completar”:" bullet below it.
#   - A long list of new bullet points (new activities: favorito/completar,
#     login, mover boton "Nueva Tarea", Modal Eliminar) is appended at the
#     end of the list, re-using the existing "Prrafodelista" numbered list
#     (numId 1), with ilvl 0 for top-level activity headers and ilvl 1 for
#     their sub-bullets.
#
# Strategy: remove the old "_GoBack" bookmark (it will be re-created, in
# its new location, as part of the inserted XML below), then insert the
# exact OOXML for the new paragraphs right before the end of the document
# body (i.e. right after the last existing paragraph).

$d = $word.ActiveDocument

# 1) The "_GoBack" bookmark currently sits at the end of the last paragraph
#    ("Si agregas nuevas funcionalidades donde van arriba o abajo?"). The
#    edit moves it down to the end of the next bullet ("favorito/completar"),
#    so drop it here; it gets re-inserted at its new spot below.
if ($d.Bookmarks.Count -gt 0) {
    try {
        $goBack = $d.Bookmarks.Item("_GoBack")
        $goBack.Delete()
    } catch {
    }
}

# 2) Append the new bullet points (as literal WordprocessingML) right at
#    the end of the document body.
$newParagraphsXml = '<w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Actividad “favorito/completar”:</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Actividad “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>login</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>”:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Agregamos “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>loginController</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">” con su método </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>index</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> en el cual </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>esta</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> la lógica de validación de usuario y contraseña con un </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>array</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> llamado post que trae esos datos de la vista del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>login</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> que </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>esta</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> en /</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>views</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>loginview.php</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Agregamos “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>loginview</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">” donde encontramos la vista del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>login</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> que manda los datos al controlador por medio de un </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>post ?</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>¿.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">Agregamos el apartado de validación de usuario correcto al </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>header</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> con su lógica de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>redireccionamiento</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> en el caso de no ser usuario correcto.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve">Agregamos un método de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>logout</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>loginController</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> que se llama desde el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>nav</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>header</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t xml:space="preserve"> el cual cierra la sesión.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Actividad mover el botón “Nueva Tarea”:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Eliminamos la tarjeta de nueva tarea del lateral.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Colocamos el botón nueva tarea en el mismo div donde se encontraba el titulo con ciertas clases para su correcta visualización y alineación.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Actividad “Modal Eliminar”:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="es-AR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-AR"/></w:rPr><w:t>Agregamos ciertos atributos al botón/link de eliminar para abrir el modal y enviar el id el cual capturamos con JavaScript para luego hacer la eliminación de la tarea.</w:t></w:r></w:p>'

$insertionPoint = $d.Content
$insertionPoint.Collapse(0)
[void]$insertionPoint.InsertXML($newParagraphsXml)
